# "adding averages and more checks"
#
# Workbook-wide header styling touch-up (shared style table affects both the
# "Training Dashboard" and "Exam Dashboard" sheets, since they reuse the same
# title/header styles) plus a couple of "Exam Dashboard" specific fixes:
#   - Header band (row 2) and report title (row 1) text becomes bold white on
#     the dark-blue band, consistent across both sheets.
#   - The oversized "COMMENTS" column on the Exam Dashboard is narrowed to
#     line up with the other data columns.
#   - The stale "redo your exam" remarks are replaced now that the exam
#     dates have been re-checked and confirmed valid.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# Header bands -> white font color (bold/fill/border stay as-is).
$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215

# Report titles -> same bold white look, regular (non-enlarged) size.
$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Color = 16777215
$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215

# Narrow the COMMENTS column (E) on the Exam Dashboard so it matches the
# width already used by the other data columns (B and C).
$ws2.Columns.Item(5).ColumnWidth = $ws2.Range("B1").ColumnWidth

# The exam dates have been reviewed and are no longer outdated.
$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"
